$wb = $excel.ActiveWorkbook

# Work on the "iOS" worksheet (sheet2) - add a new row for a Slack chat entry
$ws = $wb.Worksheets.Item("iOS")

$ws.Range("A3").Value = "Slack Chats - iOS.xml"
$ws.Range("B3").Value = "Chats"

# Select the newly added cell and make this sheet the active one
$ws.Range("B3").Select()
$ws.Activate()
